# BatxHead.xlsx — reword the two CdCode note cells on the "DBD" sheet
# (period-style "0.xxx" markers become colon-style "0:xxx"), and move the
# active selection from G17 to G16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G14 (row for BatxExeCode/作業狀態): "CdCode:BatchStatus" -> "CdCode.BatchStatus",
# and each numbered line switches from "N.xxx" to "N:xxx".
$ws.Range("G14").Value = "CdCode.BatchStatus`n0:待檢核`n1:檢核有誤`n2:檢核正常`n3:入帳未完`n4:入帳完成`n8:已刪除"

# G15 (row for BatxStsCode/整批作業狀態): "0.正常/1.整批處理中" -> "0:正常/1:整批處理中".
$ws.Range("G15").Value = "0:正常`n1:整批處理中"

# Move the saved selection/active cell from G17 to G16.
$ws.Range("G16").Select()
